$wb = $excel.ActiveWorkbook

# 1. Rename header in "Weekly Quantity" sheet (sheet1)
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2. Rename header in "Monthly Trend" sheet (sheet2)
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3. Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# 4. Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the header style (bold, border, centered) from the Weekly Quantity sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. Data rows (2-52)
$wsForecast.Cells.Item(2, 1).Value = 45025.99999999999
$wsForecast.Cells.Item(2, 2).Value = 31
$wsForecast.Cells.Item(2, 3).Value = -28.69804812503251
$wsForecast.Cells.Item(2, 4).Value = 89.48708512933696
$wsForecast.Cells.Item(3, 1).Value = 45039.99999999999
$wsForecast.Cells.Item(3, 2).Value = 32
$wsForecast.Cells.Item(3, 3).Value = -27.37461735754702
$wsForecast.Cells.Item(3, 4).Value = 92.18094887475286
$wsForecast.Cells.Item(4, 1).Value = 45046.99999999999
$wsForecast.Cells.Item(4, 2).Value = 32
$wsForecast.Cells.Item(4, 3).Value = -28.92963101117558
$wsForecast.Cells.Item(4, 4).Value = 94.35400167689906
$wsForecast.Cells.Item(5, 1).Value = 45053.99999999999
$wsForecast.Cells.Item(5, 2).Value = 32
$wsForecast.Cells.Item(5, 3).Value = -23.8414799834698
$wsForecast.Cells.Item(5, 4).Value = 92.53113382079488
$wsForecast.Cells.Item(6, 1).Value = 45060.99999999999
$wsForecast.Cells.Item(6, 2).Value = 32
$wsForecast.Cells.Item(6, 3).Value = -23.42589090152743
$wsForecast.Cells.Item(6, 4).Value = 89.65191501211004
$wsForecast.Cells.Item(7, 1).Value = 45067.99999999999
$wsForecast.Cells.Item(7, 2).Value = 32
$wsForecast.Cells.Item(7, 3).Value = -26.67330435611733
$wsForecast.Cells.Item(7, 4).Value = 88.9888547849998
$wsForecast.Cells.Item(8, 1).Value = 45074.99999999999
$wsForecast.Cells.Item(8, 2).Value = 32
$wsForecast.Cells.Item(8, 3).Value = -26.5660967907991
$wsForecast.Cells.Item(8, 4).Value = 90.12320653476678
$wsForecast.Cells.Item(9, 1).Value = 45081.99999999999
$wsForecast.Cells.Item(9, 2).Value = 32
$wsForecast.Cells.Item(9, 3).Value = -25.67632783070826
$wsForecast.Cells.Item(9, 4).Value = 88.13746332910591
$wsForecast.Cells.Item(10, 1).Value = 45088.99999999999
$wsForecast.Cells.Item(10, 2).Value = 33
$wsForecast.Cells.Item(10, 3).Value = -24.97179648945499
$wsForecast.Cells.Item(10, 4).Value = 88.15576620328223
$wsForecast.Cells.Item(11, 1).Value = 45102.99999999999
$wsForecast.Cells.Item(11, 2).Value = 33
$wsForecast.Cells.Item(11, 3).Value = -29.98334151946988
$wsForecast.Cells.Item(11, 4).Value = 90.30658744717998
$wsForecast.Cells.Item(12, 1).Value = 45116.99999999999
$wsForecast.Cells.Item(12, 2).Value = 33
$wsForecast.Cells.Item(12, 3).Value = -24.6673367056483
$wsForecast.Cells.Item(12, 4).Value = 93.49758403100985
$wsForecast.Cells.Item(13, 1).Value = 45130.99999999999
$wsForecast.Cells.Item(13, 2).Value = 33
$wsForecast.Cells.Item(13, 3).Value = -23.71679524715764
$wsForecast.Cells.Item(13, 4).Value = 89.09743655585682
$wsForecast.Cells.Item(14, 1).Value = 45137.99999999999
$wsForecast.Cells.Item(14, 2).Value = 33
$wsForecast.Cells.Item(14, 3).Value = -23.78301863047033
$wsForecast.Cells.Item(14, 4).Value = 90.05036396650544
$wsForecast.Cells.Item(15, 1).Value = 45179.99999999999
$wsForecast.Cells.Item(15, 2).Value = 34
$wsForecast.Cells.Item(15, 3).Value = -23.52922017532621
$wsForecast.Cells.Item(15, 4).Value = 94.50222873612206
$wsForecast.Cells.Item(16, 1).Value = 45186.99999999999
$wsForecast.Cells.Item(16, 2).Value = 34
$wsForecast.Cells.Item(16, 3).Value = -29.85887421618303
$wsForecast.Cells.Item(16, 4).Value = 90.82583334650619
$wsForecast.Cells.Item(17, 1).Value = 45200.99999999999
$wsForecast.Cells.Item(17, 2).Value = 34
$wsForecast.Cells.Item(17, 3).Value = -22.14736684725549
$wsForecast.Cells.Item(17, 4).Value = 90.64049977222193
$wsForecast.Cells.Item(18, 1).Value = 45214.99999999999
$wsForecast.Cells.Item(18, 2).Value = 35
$wsForecast.Cells.Item(18, 3).Value = -22.39333600982707
$wsForecast.Cells.Item(18, 4).Value = 91.21921187507151
$wsForecast.Cells.Item(19, 1).Value = 45228.99999999999
$wsForecast.Cells.Item(19, 2).Value = 35
$wsForecast.Cells.Item(19, 3).Value = -22.74131139532692
$wsForecast.Cells.Item(19, 4).Value = 91.10732838038103
$wsForecast.Cells.Item(20, 1).Value = 45242.99999999999
$wsForecast.Cells.Item(20, 2).Value = 35
$wsForecast.Cells.Item(20, 3).Value = -25.87159464078011
$wsForecast.Cells.Item(20, 4).Value = 91.24079465194369
$wsForecast.Cells.Item(21, 1).Value = 45256.99999999999
$wsForecast.Cells.Item(21, 2).Value = 35
$wsForecast.Cells.Item(21, 3).Value = -17.1422285207925
$wsForecast.Cells.Item(21, 4).Value = 96.65185881072097
$wsForecast.Cells.Item(22, 1).Value = 45277.99999999999
$wsForecast.Cells.Item(22, 2).Value = 36
$wsForecast.Cells.Item(22, 3).Value = -22.05269661479415
$wsForecast.Cells.Item(22, 4).Value = 91.2401076245325
$wsForecast.Cells.Item(23, 1).Value = 45298.99999999999
$wsForecast.Cells.Item(23, 2).Value = 36
$wsForecast.Cells.Item(23, 3).Value = -21.86130228104754
$wsForecast.Cells.Item(23, 4).Value = 96.59209267325055
$wsForecast.Cells.Item(24, 1).Value = 45305.99999999999
$wsForecast.Cells.Item(24, 2).Value = 36
$wsForecast.Cells.Item(24, 3).Value = -22.48072973505495
$wsForecast.Cells.Item(24, 4).Value = 95.21660460275989
$wsForecast.Cells.Item(25, 1).Value = 45312.99999999999
$wsForecast.Cells.Item(25, 2).Value = 36
$wsForecast.Cells.Item(25, 3).Value = -18.95257634338222
$wsForecast.Cells.Item(25, 4).Value = 97.00834512599715
$wsForecast.Cells.Item(26, 1).Value = 45347.99999999999
$wsForecast.Cells.Item(26, 2).Value = 37
$wsForecast.Cells.Item(26, 3).Value = -19.50085730519889
$wsForecast.Cells.Item(26, 4).Value = 95.83698369016508
$wsForecast.Cells.Item(27, 1).Value = 45354.99999999999
$wsForecast.Cells.Item(27, 2).Value = 37
$wsForecast.Cells.Item(27, 3).Value = -22.61283525798664
$wsForecast.Cells.Item(27, 4).Value = 93.83399242130008
$wsForecast.Cells.Item(28, 1).Value = 45361.99999999999
$wsForecast.Cells.Item(28, 2).Value = 37
$wsForecast.Cells.Item(28, 3).Value = -20.94083941740336
$wsForecast.Cells.Item(28, 4).Value = 97.87929102200447
$wsForecast.Cells.Item(29, 1).Value = 45382.99999999999
$wsForecast.Cells.Item(29, 2).Value = 38
$wsForecast.Cells.Item(29, 3).Value = -17.87222012493175
$wsForecast.Cells.Item(29, 4).Value = 101.0021156611109
$wsForecast.Cells.Item(30, 1).Value = 45396.99999999999
$wsForecast.Cells.Item(30, 2).Value = 38
$wsForecast.Cells.Item(30, 3).Value = -20.58330851972489
$wsForecast.Cells.Item(30, 4).Value = 100.033821536561
$wsForecast.Cells.Item(31, 1).Value = 45410.99999999999
$wsForecast.Cells.Item(31, 2).Value = 38
$wsForecast.Cells.Item(31, 3).Value = -20.68896095762401
$wsForecast.Cells.Item(31, 4).Value = 96.46891663273735
$wsForecast.Cells.Item(32, 1).Value = 45424.99999999999
$wsForecast.Cells.Item(32, 2).Value = 38
$wsForecast.Cells.Item(32, 3).Value = -19.3162474958668
$wsForecast.Cells.Item(32, 4).Value = 95.32477243268411
$wsForecast.Cells.Item(33, 1).Value = 45438.99999999999
$wsForecast.Cells.Item(33, 2).Value = 39
$wsForecast.Cells.Item(33, 3).Value = -19.5336826016801
$wsForecast.Cells.Item(33, 4).Value = 94.51564494216981
$wsForecast.Cells.Item(34, 1).Value = 45459.99999999999
$wsForecast.Cells.Item(34, 2).Value = 39
$wsForecast.Cells.Item(34, 3).Value = -18.67295684265401
$wsForecast.Cells.Item(34, 4).Value = 98.35827833071785
$wsForecast.Cells.Item(35, 1).Value = 45473.99999999999
$wsForecast.Cells.Item(35, 2).Value = 39
$wsForecast.Cells.Item(35, 3).Value = -20.70285394614969
$wsForecast.Cells.Item(35, 4).Value = 97.8573970913896
$wsForecast.Cells.Item(36, 1).Value = 45480.99999999999
$wsForecast.Cells.Item(36, 2).Value = 39
$wsForecast.Cells.Item(36, 3).Value = -21.89367816577182
$wsForecast.Cells.Item(36, 4).Value = 96.16341379293951
$wsForecast.Cells.Item(37, 1).Value = 45487.99999999999
$wsForecast.Cells.Item(37, 2).Value = 40
$wsForecast.Cells.Item(37, 3).Value = -16.5215755146254
$wsForecast.Cells.Item(37, 4).Value = 99.04092425055997
$wsForecast.Cells.Item(38, 1).Value = 45494.99999999999
$wsForecast.Cells.Item(38, 2).Value = 40
$wsForecast.Cells.Item(38, 3).Value = -22.652808968747
$wsForecast.Cells.Item(38, 4).Value = 100.1773767591033
$wsForecast.Cells.Item(39, 1).Value = 45501.99999999999
$wsForecast.Cells.Item(39, 2).Value = 40
$wsForecast.Cells.Item(39, 3).Value = -18.82176294772192
$wsForecast.Cells.Item(39, 4).Value = 103.2376352662957
$wsForecast.Cells.Item(40, 1).Value = 45508.99999999999
$wsForecast.Cells.Item(40, 2).Value = 40
$wsForecast.Cells.Item(40, 3).Value = -18.13678569054331
$wsForecast.Cells.Item(40, 4).Value = 97.97803559379845
$wsForecast.Cells.Item(41, 1).Value = 45515.99999999999
$wsForecast.Cells.Item(41, 2).Value = 40
$wsForecast.Cells.Item(41, 3).Value = -19.81686385775336
$wsForecast.Cells.Item(41, 4).Value = 96.47100310555579
$wsForecast.Cells.Item(42, 1).Value = 45522.99999999999
$wsForecast.Cells.Item(42, 2).Value = 40
$wsForecast.Cells.Item(42, 3).Value = -17.42952512100729
$wsForecast.Cells.Item(42, 4).Value = 94.93437580397764
$wsForecast.Cells.Item(43, 1).Value = 45529.99999999999
$wsForecast.Cells.Item(43, 2).Value = 40
$wsForecast.Cells.Item(43, 3).Value = -18.91041065155362
$wsForecast.Cells.Item(43, 4).Value = 95.5323792372
$wsForecast.Cells.Item(44, 1).Value = 45536.99999999999
$wsForecast.Cells.Item(44, 2).Value = 40
$wsForecast.Cells.Item(44, 3).Value = -18.29625800176168
$wsForecast.Cells.Item(44, 4).Value = 99.03789232514532
$wsForecast.Cells.Item(45, 1).Value = 45543.99999999999
$wsForecast.Cells.Item(45, 2).Value = 40
$wsForecast.Cells.Item(45, 3).Value = -21.74185025410653
$wsForecast.Cells.Item(45, 4).Value = 98.97331952887018
$wsForecast.Cells.Item(46, 1).Value = 45550.99999999999
$wsForecast.Cells.Item(46, 2).Value = 41
$wsForecast.Cells.Item(46, 3).Value = -17.63222850196815
$wsForecast.Cells.Item(46, 4).Value = 96.04121145712679
$wsForecast.Cells.Item(47, 1).Value = 45557.99999999999
$wsForecast.Cells.Item(47, 2).Value = 41
$wsForecast.Cells.Item(47, 3).Value = -18.37968795595319
$wsForecast.Cells.Item(47, 4).Value = 96.2299109753081
$wsForecast.Cells.Item(48, 1).Value = 45564.99999999999
$wsForecast.Cells.Item(48, 2).Value = 41
$wsForecast.Cells.Item(48, 3).Value = -18.65795625599603
$wsForecast.Cells.Item(48, 4).Value = 98.56006602444221
$wsForecast.Cells.Item(49, 1).Value = 45571.99999999999
$wsForecast.Cells.Item(49, 2).Value = 41
$wsForecast.Cells.Item(49, 3).Value = -17.98444314072713
$wsForecast.Cells.Item(49, 4).Value = 93.91294888874583
$wsForecast.Cells.Item(50, 1).Value = 45578.99999999999
$wsForecast.Cells.Item(50, 2).Value = 41
$wsForecast.Cells.Item(50, 3).Value = -15.2430893725137
$wsForecast.Cells.Item(50, 4).Value = 104.5736850519081
$wsForecast.Cells.Item(51, 1).Value = 45585.99999999999
$wsForecast.Cells.Item(51, 2).Value = 41
$wsForecast.Cells.Item(51, 3).Value = -19.02389356438165
$wsForecast.Cells.Item(51, 4).Value = 100.3984542671511
$wsForecast.Cells.Item(52, 1).Value = 45592.99999999999
$wsForecast.Cells.Item(52, 2).Value = 41
$wsForecast.Cells.Item(52, 3).Value = -14.39943942251821
$wsForecast.Cells.Item(52, 4).Value = 101.4809700536036

# Copy the date-format style from column A of the Weekly Quantity sheet to the ds column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A52").PasteSpecial(-4122)
$excel.CutCopyMode = $false

